# Update loading percent results for case with 380 kV
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 14.22178750299507
$ws.Cells.Item(2, 3).Value = 9.39369637722708
$ws.Cells.Item(2, 4).Value = 9.703837241471204
$ws.Cells.Item(2, 5).Value = 13.86929718835303
$ws.Cells.Item(2, 6).Value = 30.11101205632314
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 9).Value = 19.67133455059014
$ws.Cells.Item(2, 10).Value = 9.798452562305377
$ws.Cells.Item(2, 14).Value = 16.57561220254406
$ws.Cells.Item(2, 15).Value = 22.23624642048586
$ws.Cells.Item(3, 2).Value = 13.66403666633301
$ws.Cells.Item(3, 3).Value = 8.894706724168175
$ws.Cells.Item(3, 4).Value = 9.659539720370843
$ws.Cells.Item(3, 5).Value = 13.81211111198105
$ws.Cells.Item(3, 6).Value = 30.09216981917425
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 9).Value = 19.77229499921408
$ws.Cells.Item(3, 10).Value = 9.8039936253913
$ws.Cells.Item(3, 14).Value = 16.61480622921754
$ws.Cells.Item(3, 15).Value = 22.27421354860361
$ws.Cells.Item(4, 2).Value = 13.31105359881824
$ws.Cells.Item(4, 3).Value = 8.573939683486827
$ws.Cells.Item(4, 4).Value = 9.633791041158243
$ws.Cells.Item(4, 5).Value = 13.77976551609274
$ws.Cells.Item(4, 6).Value = 30.08933451730041
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 9).Value = 19.8390040439844
$ws.Cells.Item(4, 10).Value = 9.809027471025177
$ws.Cells.Item(4, 14).Value = 16.64065590125057
$ws.Cells.Item(4, 15).Value = 22.30359263973622
$ws.Cells.Item(5, 2).Value = 13.16477447573105
$ws.Cells.Item(5, 3).Value = 8.439734575707948
$ws.Cells.Item(5, 4).Value = 9.623670518990689
$ws.Cells.Item(5, 5).Value = 13.76729008282449
$ws.Cells.Item(5, 6).Value = 30.09037584019267
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 9).Value = 19.86737241629787
$ws.Cells.Item(5, 10).Value = 9.811489339460657
$ws.Cells.Item(5, 14).Value = 16.65163923458547
$ws.Cells.Item(5, 15).Value = 22.31708586296834
$ws.Cells.Item(6, 2).Value = 13.14034444435391
$ws.Cells.Item(6, 3).Value = 8.417243336494391
$ws.Cells.Item(6, 4).Value = 9.622012712074399
$ws.Cells.Item(6, 5).Value = 13.76526143887645
$ws.Cells.Item(6, 6).Value = 30.09068140754522
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 9).Value = 19.87215438875312
$ws.Cells.Item(6, 10).Value = 9.811922931155481
$ws.Cells.Item(6, 14).Value = 16.65349017066815
$ws.Cells.Item(6, 15).Value = 22.31941812057849
$ws.Cells.Item(7, 2).Value = 13.30909039399474
$ws.Cells.Item(7, 3).Value = 8.572143691079996
$ws.Cells.Item(7, 4).Value = 9.633653035206899
$ws.Cells.Item(7, 5).Value = 13.77959439841244
$ws.Cells.Item(7, 6).Value = 30.0893396671761
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 9).Value = 19.8393818398094
$ws.Cells.Item(7, 10).Value = 9.809059010158309
$ws.Cells.Item(7, 14).Value = 16.64080220588432
$ws.Cells.Item(7, 15).Value = 22.30376846215606
$ws.Cells.Item(8, 2).Value = 14.03177754272028
$ws.Cells.Item(8, 3).Value = 9.224699694391637
$ws.Cells.Item(8, 4).Value = 9.688267172098959
$ws.Cells.Item(8, 5).Value = 13.84901136460695
$ws.Cells.Item(8, 6).Value = 30.10270338333554
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 9).Value = 19.70516465932628
$ws.Cells.Item(8, 10).Value = 9.800024666501821
$ws.Cells.Item(8, 14).Value = 16.58875636316404
$ws.Cells.Item(8, 15).Value = 22.24807565219315
$ws.Cells.Item(9, 2).Value = 15.35742809848736
$ws.Cells.Item(9, 3).Value = 10.38572983715067
$ws.Cells.Item(9, 4).Value = 9.806518158539882
$ws.Cells.Item(9, 5).Value = 14.00661780036621
$ws.Cells.Item(9, 6).Value = 30.19810684302134
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 9).Value = 19.47953020030432
$ws.Cells.Item(9, 10).Value = 9.79523935019774
$ws.Cells.Item(9, 14).Value = 16.50082290419505
$ws.Cells.Item(9, 15).Value = 22.18718336166931
$ws.Cells.Item(10, 2).Value = 16.26633074300051
$ws.Cells.Item(10, 3).Value = 11.16169228561528
$ws.Cells.Item(10, 4).Value = 9.899696528500373
$ws.Cells.Item(10, 5).Value = 14.13481746408114
$ws.Cells.Item(10, 6).Value = 30.31012993904219
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 9).Value = 19.33682660710959
$ws.Cells.Item(10, 10).Value = 9.799580493783678
$ws.Cells.Item(10, 14).Value = 16.44479100444371
$ws.Cells.Item(10, 15).Value = 22.17211329713646
$ws.Cells.Item(11, 2).Value = 16.66413035629505
$ws.Cells.Item(11, 3).Value = 11.49726837044938
$ws.Cells.Item(11, 4).Value = 9.943339619410892
$ws.Cells.Item(11, 5).Value = 14.1956702546349
$ws.Cells.Item(11, 6).Value = 30.37010287608231
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 9).Value = 19.27695396906568
$ws.Cells.Item(11, 10).Value = 9.803253196481384
$ws.Cells.Item(11, 14).Value = 16.42115394043351
$ws.Cells.Item(11, 15).Value = 22.17173033711716
$ws.Cells.Item(12, 2).Value = 16.81240739323342
$ws.Cells.Item(12, 3).Value = 11.62179498719444
$ws.Cells.Item(12, 4).Value = 9.960036580241081
$ws.Cells.Item(12, 5).Value = 14.21906326675251
$ws.Cells.Item(12, 6).Value = 30.39409804200044
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 9).Value = 19.25501035278481
$ws.Cells.Item(12, 10).Value = 9.804887145941425
$ws.Cells.Item(12, 14).Value = 16.4124689268118
$ws.Cells.Item(12, 15).Value = 22.17251733069078
$ws.Cells.Item(13, 2).Value = 16.78057985393075
$ws.Cells.Item(13, 3).Value = 11.59508992327099
$ws.Cells.Item(13, 4).Value = 9.956433185167219
$ws.Cells.Item(13, 5).Value = 14.21400986972647
$ws.Cells.Item(13, 6).Value = 30.38887331656151
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 9).Value = 19.25970383035609
$ws.Cells.Item(13, 10).Value = 9.804524447450257
$ws.Cells.Item(13, 14).Value = 16.41432758491269
$ws.Cells.Item(13, 15).Value = 22.1723063738193
$ws.Cells.Item(14, 2).Value = 16.67637708329384
$ws.Cells.Item(14, 3).Value = 11.50756455442098
$ws.Cells.Item(14, 4).Value = 9.944709935969517
$ws.Cells.Item(14, 5).Value = 14.19758789786512
$ws.Cells.Item(14, 6).Value = 30.37205129698316
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 9).Value = 19.27513402680545
$ws.Cells.Item(14, 10).Value = 9.803382754978694
$ws.Cells.Item(14, 14).Value = 16.42043409434974
$ws.Cells.Item(14, 15).Value = 22.17177640003412
$ws.Cells.Item(15, 2).Value = 16.61223944208493
$ws.Cells.Item(15, 3).Value = 11.45361964621836
$ws.Cells.Item(15, 4).Value = 9.937550970903807
$ws.Cells.Item(15, 5).Value = 14.18757400959682
$ws.Cells.Item(15, 6).Value = 30.36191426282698
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 9).Value = 19.28468049168449
$ws.Cells.Item(15, 10).Value = 9.802715072665505
$ws.Cells.Item(15, 14).Value = 16.42420911185271
$ws.Cells.Item(15, 15).Value = 22.17157317619421
$ws.Cells.Item(16, 2).Value = 16.24000876108173
$ws.Cells.Item(16, 3).Value = 11.13940715544861
$ws.Cells.Item(16, 4).Value = 9.89686874427135
$ws.Cells.Item(16, 5).Value = 14.13089031381781
$ws.Cells.Item(16, 6).Value = 30.30639100576419
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 9).Value = 19.3408412232166
$ws.Cells.Item(16, 10).Value = 9.799374559633005
$ws.Cells.Item(16, 14).Value = 16.44637296825818
$ws.Cells.Item(16, 15).Value = 22.1722686938948
$ws.Cells.Item(17, 2).Value = 16.00756301475243
$ws.Cells.Item(17, 3).Value = 10.94215308315854
$ws.Cells.Item(17, 4).Value = 9.872225904028003
$ws.Cells.Item(17, 5).Value = 14.09675525868707
$ws.Cells.Item(17, 6).Value = 30.27463028012304
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 9).Value = 19.3765883895725
$ws.Cells.Item(17, 10).Value = 9.797759483925255
$ws.Cells.Item(17, 14).Value = 16.4604437869993
$ws.Cells.Item(17, 15).Value = 22.17435427923074
$ws.Cells.Item(18, 2).Value = 15.87239736474448
$ws.Cells.Item(18, 3).Value = 10.82706075515996
$ws.Cells.Item(18, 4).Value = 9.85817082993441
$ws.Cells.Item(18, 5).Value = 14.07736114111405
$ws.Cells.Item(18, 6).Value = 30.25721150532636
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 9).Value = 19.39762375882542
$ws.Cells.Item(18, 10).Value = 9.796990416609429
$ws.Cells.Item(18, 14).Value = 16.46871130996559
$ws.Cells.Item(18, 15).Value = 22.17616304516275
$ws.Cells.Item(19, 2).Value = 15.82638385594282
$ws.Cells.Item(19, 3).Value = 10.78781283687904
$ws.Cells.Item(19, 4).Value = 9.85343274317141
$ws.Cells.Item(19, 5).Value = 14.07083620188029
$ws.Cells.Item(19, 6).Value = 30.2514599651899
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 9).Value = 19.40482736995183
$ws.Cells.Item(19, 10).Value = 9.796757510904005
$ws.Cells.Item(19, 14).Value = 16.4715405133506
$ws.Cells.Item(19, 15).Value = 22.17688004035956
$ws.Cells.Item(20, 2).Value = 16.03246013238663
$ws.Cells.Item(20, 3).Value = 10.96332089404367
$ws.Cells.Item(20, 4).Value = 9.874836951210106
$ws.Cells.Item(20, 5).Value = 14.10036431692925
$ws.Cells.Item(20, 6).Value = 30.2779234554144
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 9).Value = 19.37273390359694
$ws.Cells.Item(20, 10).Value = 9.797914870988441
$ws.Cells.Item(20, 14).Value = 16.4589278823909
$ws.Cells.Item(20, 15).Value = 22.17406920429874
$ws.Cells.Item(21, 2).Value = 16.7070488298273
$ws.Cells.Item(21, 3).Value = 11.53334236010009
$ws.Cells.Item(21, 4).Value = 9.948148800849326
$ws.Cells.Item(21, 5).Value = 14.20240206964578
$ws.Cells.Item(21, 6).Value = 30.37695756497387
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 9).Value = 19.27058199476645
$ws.Cells.Item(21, 10).Value = 9.803711506497299
$ws.Cells.Item(21, 14).Value = 16.41863325390795
$ws.Cells.Item(21, 15).Value = 22.17190676523548
$ws.Cells.Item(22, 2).Value = 17.13412968299703
$ws.Cells.Item(22, 3).Value = 11.89101370767598
$ws.Cells.Item(22, 4).Value = 9.997049832952735
$ws.Cells.Item(22, 5).Value = 14.27111872608314
$ws.Cells.Item(22, 6).Value = 30.44916379081098
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 9).Value = 19.20806980471752
$ws.Cells.Item(22, 10).Value = 9.808916807395217
$ws.Cells.Item(22, 14).Value = 16.39384761975928
$ws.Cells.Item(22, 15).Value = 22.17592601907452
$ws.Cells.Item(23, 2).Value = 16.90748343254125
$ws.Cells.Item(23, 3).Value = 11.70149088636061
$ws.Cells.Item(23, 4).Value = 9.97086354801713
$ws.Cells.Item(23, 5).Value = 14.23426280502721
$ws.Cells.Item(23, 6).Value = 30.40994563191617
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 9).Value = 19.24104368024516
$ws.Cells.Item(23, 10).Value = 9.806009357464816
$ws.Cells.Item(23, 14).Value = 16.40693458666685
$ws.Cells.Item(23, 15).Value = 22.1732835661038
$ws.Cells.Item(24, 2).Value = 16.0212089011229
$ws.Cells.Item(24, 3).Value = 10.95375618206884
$ws.Cells.Item(24, 4).Value = 9.873656145821391
$ws.Cells.Item(24, 5).Value = 14.09873194234462
$ws.Cells.Item(24, 6).Value = 30.27643199067597
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 9).Value = 19.37447500978024
$ws.Cells.Item(24, 10).Value = 9.797844123719148
$ws.Cells.Item(24, 14).Value = 16.45961266828392
$ws.Cells.Item(24, 15).Value = 22.17419618745662
$ws.Cells.Item(25, 2).Value = 15.00963879198647
$ws.Cells.Item(25, 3).Value = 10.08489064325125
$ws.Cells.Item(25, 4).Value = 9.773384200784182
$ws.Cells.Item(25, 5).Value = 13.96174946392132
$ws.Cells.Item(25, 6).Value = 30.16490827828055
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 9).Value = 19.53653142203374
$ws.Cells.Item(25, 10).Value = 9.795151714684291
$ws.Cells.Item(25, 14).Value = 16.52310291909228
$ws.Cells.Item(25, 15).Value = 22.19845889362177
